$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same style as the other header cells (e.g. H1) to I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

# Data values for columns I (I0) and J (IF), rows 2-23
$iValues = @(8,7,4,9,7,9,8,7,8,8,4,8,4,6,8,4,8,6,7,9,6,7)
$jValues = @(9,7,5,9,8,9,9,7,8,8,4,8,5,6,8,4,8,6,7,9,6,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
